$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.847.14'
$ws.Range('E2').Value = '  +4.79%  '
$ws.Range('D3').Value = '1.867.28'
$ws.Range('E3').Value = '  +3.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '272.66'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9995'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5294'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +5.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3383'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06814'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.89'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7936'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07748'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.71%  '
$ws.Range('D13').Value = '1.847.53'
$ws.Range('E13').Value = '  +1.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '90.31'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.131'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9987'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.40'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008003'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9998'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').Value = '26.880.40'
$ws.Range('E20').Value = '  +4.70%  '
$ws.Range('D21').Value = '2.115.15'
$ws.Range('E21').Value = '  +2.83%  '
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.937'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.066'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.386'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.79'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.656'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '112.64'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.326'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.303'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08855'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04946'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.162'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7278'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.870'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.192'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.66%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01846'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.40%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.309'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5089'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '116.09'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.12%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9347'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.141'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.016'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9991'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4421'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1327'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.325'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.13'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05937'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.468'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.79%  '
